# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Swap the "Periodo Mora" / "Valor Mora" figures between the first and third
# data rows of the account-statement table (row 16 and row 18), leaving the
# middle row (17) untouched.
#   Row 16: Periodo Mora 2407 -> 2402 , Valor Mora 677761 -> 638622
#   Row 18: Periodo Mora 2402 -> 2407 , Valor Mora 638622 -> 677761

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these as text (quoted) so the period codes remain stored as strings,
# matching the existing "Periodo Mora" column formatting.
$ws.Range("E16").Value = "2402"
$ws.Range("E18").Value = "2407"

$ws.Range("F16").Value = 638622
$ws.Range("F18").Value = 677761
